$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert two new line-item rows after row 14 ("mouse" and "monitor"), matching
#     the existing "keyboard" row's layout (if not empty -> add another PR line) ---
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(16).Insert()

# Re-apply row 14's exact formatting onto the two new rows.
$ws.Range("A14:K14").Copy()
$ws.Range("A15:K15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A14:K14").Copy()
$ws.Range("A16:K16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the merged ranges (Description spans E:H, Date Needed spans J:K).
$ws.Range("E15:H15").Merge()
$ws.Range("J15:K15").Merge()
$ws.Range("E16:H16").Merge()
$ws.Range("J16:K16").Merge()

# --- Quantities now that there are three line items instead of one ---
$ws.Range("B14").Value = 1

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "pc/s"
$ws.Range("E15").Value = "mouse"

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "pc/s"
$ws.Range("E16").Value = "monitor"

# --- Leave the cursor where Excel would land after editing this table ---
$ws.Range("B20").Select()
